$d = $word.ActiveDocument

# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll) -- each needle below is unique in
# the document, so "replace all" vs. "replace one" makes no practical difference here.

# 1) "POR NORESTE MIDE:" -> "POR NORTE MIDE:"
#    (drops the stray "ES" out of "NORESTE"; the run is simply re-typed in place so its
#    Arial Narrow / bold / black formatting carries over unchanged)
$d.Content.Find.Execute("POR NORESTE MIDE:", $true, $false, $false, $false, $false, $true, 1, $false, "POR NORTE MIDE:", 2)

# 2) "POR SUROESTE MIDE:" -> "POR SUR MIDE:"
$d.Content.Find.Execute("POR SUROESTE MIDE:", $true, $false, $false, $false, $false, $true, 1, $false, "POR SUR MIDE:", 2)

# 3) " QUEDARÁ FACULTADO PARA DISPONER ..." -> " QUEDARÁ FACULTAD{{SEXO_9}} PARA DISPONER ..."
#    (the trailing "O" of "FACULTADO" becomes the gender merge-field {{SEXO_9}}, matching the
#    {{SEXO_##}} placeholders used throughout the rest of the contract)
$d.Content.Find.Execute("QUEDARÁ FACULTADO PARA DISPONER", $true, $false, $false, $false, $false, $true, 1, $false, "QUEDARÁ FACULTAD{{SEXO_9}} PARA DISPONER", 2)
